$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" holds the text block with conversion rates in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.94 = 10892.86 pesos`n✅ 10892.86 pesos = 2.92 = 953.43 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas" holds the numeric rate cells N10, O10, N12, O12
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 339.999
$ws2.Range("O10").Value = 3703.56
$ws2.Range("N12").Value = 3724.52
$ws2.Range("O12").Value = 326
